$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 5 new rows at the positions where new weekly records were added
# (ascending order so each index refers to the already-shifted sheet)
$ws.Rows.Item(36).Insert()
$ws.Rows.Item(37).Insert()
$ws.Rows.Item(38).Insert()
$ws.Rows.Item(52).Insert()
$ws.Rows.Item(67).Insert()

# Row 36
$ws.Range("A36").Value = 5
$ws.Range("B36").Value = 'Macroferia Regional de Talca'
$ws.Range("C36").Value = 'Maule'
$ws.Range("D36").Value = 44435
$ws.Range("E36").Value = 7
$ws.Range("F36").Value = 'Fruta'
$ws.Range("G36").Value = 100108
$ws.Range("H36").Value = 'Tropicales y subtropicales'
$ws.Range("I36").Value = 100108002
$ws.Range("J36").Value = 'Mango'
$ws.Range("K36").Value = 'Sin especificar'
$ws.Range("L36").Value = 'Primera'
$ws.Range("M36").Value = 450
$ws.Range("N36").Value = 7500
$ws.Range("O36").Value = 8000
$ws.Range("P36").Value = 7889
$ws.Range("Q36").Value = '$/bandeja 4 kilos'
$ws.Range("R36").Value = 'Brasil'
$ws.Range("S36").Value = 1972
$ws.Range("T36").Value = 4

# Row 37
$ws.Range("A37").Value = 5
$ws.Range("B37").Value = 'Macroferia Regional de Talca'
$ws.Range("C37").Value = 'Maule'
$ws.Range("D37").Value = 44435
$ws.Range("E37").Value = 7
$ws.Range("F37").Value = 'Fruta'
$ws.Range("G37").Value = 100108
$ws.Range("H37").Value = 'Tropicales y subtropicales'
$ws.Range("I37").Value = 100108002
$ws.Range("J37").Value = 'Mango'
$ws.Range("K37").Value = 'Sin especificar'
$ws.Range("L37").Value = 'Primera'
$ws.Range("M37").Value = 120
$ws.Range("N37").Value = 8000
$ws.Range("O37").Value = 8000
$ws.Range("P37").Value = 8000
$ws.Range("Q37").Value = '$/bandeja 4 kilos'
$ws.Range("R37").Value = 'México'
$ws.Range("S37").Value = 2000
$ws.Range("T37").Value = 4

# Row 38
$ws.Range("A38").Value = 5
$ws.Range("B38").Value = 'Macroferia Regional de Talca'
$ws.Range("C38").Value = 'Maule'
$ws.Range("D38").Value = 44431
$ws.Range("E38").Value = 7
$ws.Range("F38").Value = 'Fruta'
$ws.Range("G38").Value = 100108
$ws.Range("H38").Value = 'Tropicales y subtropicales'
$ws.Range("I38").Value = 100108002
$ws.Range("J38").Value = 'Mango'
$ws.Range("K38").Value = 'Sin especificar'
$ws.Range("L38").Value = 'Primera'
$ws.Range("M38").Value = 230
$ws.Range("N38").Value = 8000
$ws.Range("O38").Value = 8000
$ws.Range("P38").Value = 8000
$ws.Range("Q38").Value = '$/bandeja 4 kilos'
$ws.Range("R38").Value = 'Brasil'
$ws.Range("S38").Value = 2000
$ws.Range("T38").Value = 4

# Row 52
$ws.Range("A52").Value = 5
$ws.Range("B52").Value = 'Macroferia Regional de Talca'
$ws.Range("C52").Value = 'Maule'
$ws.Range("D52").Value = 44433
$ws.Range("E52").Value = 7
$ws.Range("F52").Value = 'Fruta'
$ws.Range("G52").Value = 100108
$ws.Range("H52").Value = 'Tropicales y subtropicales'
$ws.Range("I52").Value = 100108002
$ws.Range("J52").Value = 'Mango'
$ws.Range("K52").Value = 'Sin especificar'
$ws.Range("L52").Value = 'Primera'
$ws.Range("M52").Value = 220
$ws.Range("N52").Value = 7500
$ws.Range("O52").Value = 8000
$ws.Range("P52").Value = 7773
$ws.Range("Q52").Value = '$/bandeja 4 kilos'
$ws.Range("R52").Value = 'Brasil'
$ws.Range("S52").Value = 1943
$ws.Range("T52").Value = 4

# Row 67
$ws.Range("A67").Value = 5
$ws.Range("B67").Value = 'Macroferia Regional de Talca'
$ws.Range("C67").Value = 'Maule'
$ws.Range("D67").Value = 44432
$ws.Range("E67").Value = 7
$ws.Range("F67").Value = 'Fruta'
$ws.Range("G67").Value = 100108
$ws.Range("H67").Value = 'Tropicales y subtropicales'
$ws.Range("I67").Value = 100108002
$ws.Range("J67").Value = 'Mango'
$ws.Range("K67").Value = 'Sin especificar'
$ws.Range("L67").Value = 'Primera'
$ws.Range("M67").Value = 120
$ws.Range("N67").Value = 8000
$ws.Range("O67").Value = 8000
$ws.Range("P67").Value = 8000
$ws.Range("Q67").Value = '$/bandeja 4 kilos'
$ws.Range("R67").Value = 'México'
$ws.Range("S67").Value = 2000
$ws.Range("T67").Value = 4
